$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value still reads as a plain number (e.g. "242.21")
# must be pre-formatted as Text, otherwise Excel auto-converts the literal
# string into a numeric value (dropping e.g. trailing zeros) on assignment -
# these Price cells are plain text in the workbook, not numbers.
$textPriceCells = @(
    "D5", "D6", "D7", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($ref in $textPriceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.883.70'
$ws.Range("E2").Value = '  +1.78%  '
$ws.Range("D3").Value = '1.730.77'
$ws.Range("E3").Value = '  +0.31%  '
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").Value = '242.21'
$ws.Range("E5").Value = '  -0.65%  '
$ws.Range("D6").Value = '0.9978'
$ws.Range("E6").Value = '  -0.25%  '
$ws.Range("D7").Value = '0.4917'
$ws.Range("E7").Value = '  +0.35%  '
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").Value = '0.06227'
$ws.Range("E9").Value = '  +0.46%  '
$ws.Range("D10").Value = '1.735.63'
$ws.Range("E10").Value = '  +0.56%  '
$ws.Range("D11").Value = '16.08'
$ws.Range("E11").Value = '  +3.54%  '
$ws.Range("D12").Value = '0.06902'
$ws.Range("E12").Value = '  -1.61%  '
$ws.Range("D13").Value = '0.6109'
$ws.Range("E13").Value = '  +1.78%  '
$ws.Range("D14").Value = '4.508'
$ws.Range("E14").Value = '  -1.66%  '
$ws.Range("D15").Value = '77.38'
$ws.Range("E15").Value = '  +0.06%  '
$ws.Range("D16").Value = '0.9981'
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("D17").Value = '26.648.78'
$ws.Range("E17").Value = '  +0.83%  '
$ws.Range("D18").Value = '0.9971'
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("D19").Value = '0.000007193'
$ws.Range("E19").Value = '  +0.60%  '
$ws.Range("D20").Value = '11.47'
$ws.Range("E20").Value = '  +1.02%  '
$ws.Range("D21").Value = '1.960.44'
$ws.Range("E21").Value = '  +0.92%  '
$ws.Range("D22").Value = '4.443'
$ws.Range("E22").Value = '  -0.68%  '
$ws.Range("D23").Value = '8.568'
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").Value = '5.129'
$ws.Range("E24").Value = '  -0.70%  '
$ws.Range("D25").Value = '138.82'
$ws.Range("E25").Value = '  +0.84%  '
$ws.Range("D26").Value = '15.32'
$ws.Range("E26").Value = '  +0.60%  '
$ws.Range("D27").Value = '1.793'
$ws.Range("E27").Value = '  +5.28%  '
$ws.Range("D28").Value = '1.381'
$ws.Range("E28").Value = '  -0.76%  '
$ws.Range("D29").Value = '106.15'
$ws.Range("E29").Value = '  -0.71%  '
$ws.Range("D30").Value = '3.947'
$ws.Range("E30").Value = '  -0.28%  '
$ws.Range("D31").Value = '0.07999'
$ws.Range("E31").Value = '  +0.62%  '
$ws.Range("D32").Value = '3.682'
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("D34").Value = '0.9968'
$ws.Range("E34").Value = '  -0.26%  '
$ws.Range("D35").Value = '2.608'
$ws.Range("E35").Value = '  +0.12%  '
$ws.Range("D36").Value = '1.012'
$ws.Range("E36").Value = '  +1.43%  '
$ws.Range("D37").Value = '0.6256'
$ws.Range("E37").Value = '  -0.30%  '
$ws.Range("D38").Value = '0.9340'
$ws.Range("E38").Value = '  +2.56%  '
$ws.Range("D39").Value = '2.061'
$ws.Range("E39").Value = '  +5.18%  '
$ws.Range("D40").Value = '2.443'
$ws.Range("E40").Value = '  +2.09%  '
$ws.Range("D41").Value = '1.000'
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("E42").Value = '  +1.30%  '
$ws.Range("D43").Value = '5.651'
$ws.Range("E43").Value = '  +3.77%  '
$ws.Range("D44").Value = '99.81'
$ws.Range("E44").Value = '  +0.36%  '
$ws.Range("D45").Value = '0.3878'
$ws.Range("E45").Value = '  +0.76%  '
$ws.Range("D46").Value = '6.953'
$ws.Range("E46").Value = '  +3.52%  '
$ws.Range("D47").Value = '0.1161'
$ws.Range("E47").Value = '  +0.24%  '
$ws.Range("D48").Value = '0.05388'
$ws.Range("E48").Value = '  +0.39%  '
$ws.Range("D49").Value = '7.977'
$ws.Range("E49").Value = '  +3.09%  '
$ws.Range("D50").Value = '30.31'
$ws.Range("E50").Value = '  +0.59%  '
$ws.Range("D51").Value = '1.244'
$ws.Range("E51").Value = '  +0.40%  '
